$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range('D2').Value = '30.325.17'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '1.868.69'
$ws.Range('E4').Value = '  +0.10%  '
Set-TextValue $ws 'D5' '235.61'
$ws.Range('E5').Value = '  +0.37%  '
Set-TextValue $ws 'D6' '1.001'
$ws.Range('E6').Value = '  +0.10%  '
Set-TextValue $ws 'D7' '0.4676'
$ws.Range('E7').Value = '  -0.52%  '
Set-TextValue $ws 'D8' '0.2845'
$ws.Range('E8').Value = '  +0.20%  '
Set-TextValue $ws 'D9' '0.06539'
$ws.Range('E9').Value = '  -1.10%  '
Set-TextValue $ws 'D10' '21.40'
$ws.Range('E10').Value = '  +5.56%  '
Set-TextValue $ws 'D11' '0.07877'
$ws.Range('E11').Value = '  +1.32%  '
Set-TextValue $ws 'D12' '97.68'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '1.872.98'
$ws.Range('E13').Value = '  -0.13%  '
Set-TextValue $ws 'D14' '5.094'
$ws.Range('E15').Value = '  +0.37%  '
Set-TextValue $ws 'D16' '276.24'
$ws.Range('E16').Value = '  -3.89%  '
$ws.Range('D17').Value = '30.319.98'
$ws.Range('E18').Value = '  +0.05%  '
Set-TextValue $ws 'D19' '12.71'
$ws.Range('E19').Value = '  +0.88%  '
Set-TextValue $ws 'D20' '5.467'
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws 'D21' '0.000007304'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.106.80'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('E23').Value = '  +0.12%  '
Set-TextValue $ws 'D24' '6.141'
$ws.Range('E24').Value = '  -0.61%  '
Set-TextValue $ws 'D25' '165.48'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('E26').Value = '  -2.71%  '
Set-TextValue $ws 'D27' '19.05'
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('E29').Value = '  +0.82%  '
Set-TextValue $ws 'D30' '0.09623'
$ws.Range('E30').Value = '  -0.53%  '
Set-TextValue $ws 'D31' '4.384'
$ws.Range('E31').Value = '  -0.22%  '
Set-TextValue $ws 'D32' '1.475'
$ws.Range('E32').Value = '  +0.61%  '
Set-TextValue $ws 'D33' '4.092'
$ws.Range('E33').Value = '  -0.75%  '
Set-TextValue $ws 'D34' '0.04699'
$ws.Range('E34').Value = '  +0.13%  '
Set-TextValue $ws 'D35' '1.126'
$ws.Range('E35').Value = '  +3.04%  '
Set-TextValue $ws 'D36' '0.7050'
$ws.Range('E36').Value = '  -0.22%  '
Set-TextValue $ws 'D37' '2.721'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('E38').Value = '  -0.84%  '
Set-TextValue $ws 'D39' '6.324'
$ws.Range('E39').Value = '  -1.83%  '
Set-TextValue $ws 'D40' '2.535'
$ws.Range('E40').Value = '  +0.39%  '
Set-TextValue $ws 'D41' '73.81'
$ws.Range('E41').Value = '  +2.55%  '
Set-TextValue $ws 'D42' '1.950'
$ws.Range('E42').Value = '  -0.26%  '
Set-TextValue $ws 'D43' '0.8496'
$ws.Range('E43').Value = '  -1.29%  '
Set-TextValue $ws 'D44' '0.4185'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('E45').Value = '  +0.10%  '
Set-TextValue $ws 'D46' '103.76'
$ws.Range('E46').Value = '  +0.71%  '
Set-TextValue $ws 'D47' '7.185'
$ws.Range('E47').Value = '  -0.40%  '
Set-TextValue $ws 'D48' '9.240'
$ws.Range('E48').Value = '  +0.43%  '
Set-TextValue $ws 'D49' '935.77'
$ws.Range('E49').Value = '  -5.46%  '
Set-TextValue $ws 'D50' '34.14'
$ws.Range('E50').Value = '  +0.27%  '
Set-TextValue $ws 'D51' '0.05633'
